$d = $word.ActiveDocument

# Remove the _GoBack bookmark from its current location (end of the
# "hbs.registerHelper(helpername,function);" paragraph).
$bm = $d.Bookmarks.Item("_GoBack")
$bm.Delete()

# The document currently ends with a trailing empty paragraph right
# before the sectPr. Insert three more paragraphs after it: one empty,
# one for the "in express module..." line, and one for the
# "app.use work in order" line.
$n = $d.Paragraphs.Count
$p = $d.Paragraphs.Item($n)
$p.Range.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$p2 = $d.Paragraphs.Item($n)
$p2.Range.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$p3 = $d.Paragraphs.Item($n)
$p3.Range.InsertParagraphAfter()

$n = $d.Paragraphs.Count
$pMiddleware = $d.Paragraphs.Item($n - 1)
$pMiddleware.Range.Text = "in express module  app.use() works as middle ware "

$pOrder = $d.Paragraphs.Item($n)
# Temporarily append a sentinel character so the insertion point for the
# new bookmark is not sitting exactly on the paragraph-mark / end-of-story
# boundary (collapsed bookmarks placed there land at the wrong offset).
$pOrder.Range.Text = "app.use work in orderX"

$bookmarkPos = $pOrder.Range.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the sentinel character now that the bookmark is anchored.
$sentinelRange = $d.Range($bookmarkPos, $bookmarkPos + 1)
$sentinelRange.Delete()
